# Update "想去人数" (F column) values on sheet "展览" and sheet "全部类型"
# to reflect freshly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 153
$ws1.Range("F4").Value  = 1784
$ws1.Range("F6").Value  = 1090
$ws1.Range("F7").Value  = 2225
$ws1.Range("F8").Value  = 2136
$ws1.Range("F10").Value = 610
$ws1.Range("F12").Value = 1683
$ws1.Range("F17").Value = 223
$ws1.Range("F18").Value = 1597
$ws1.Range("F19").Value = 642
$ws1.Range("F20").Value = 733
$ws1.Range("F21").Value = 617
$ws1.Range("F22").Value = 12292
$ws1.Range("F23").Value = 12352
$ws1.Range("F24").Value = 914
$ws1.Range("F25").Value = 705
$ws1.Range("F29").Value = 381
$ws1.Range("F30").Value = 1927
$ws1.Range("F33").Value = 204

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value  = 153
$ws4.Range("F5").Value  = 1784
$ws4.Range("F7").Value  = 1090
$ws4.Range("F8").Value  = 2225
$ws4.Range("F9").Value  = 2136
$ws4.Range("F11").Value = 610
$ws4.Range("F13").Value = 1683
$ws4.Range("F21").Value = 223
$ws4.Range("F22").Value = 1598
$ws4.Range("F23").Value = 642
$ws4.Range("F24").Value = 733
$ws4.Range("F25").Value = 617
$ws4.Range("F26").Value = 12292
$ws4.Range("F27").Value = 12352
$ws4.Range("F28").Value = 914
$ws4.Range("F29").Value = 705
$ws4.Range("F31").Value = 43
$ws4.Range("F33").Value = 381
$ws4.Range("F34").Value = 1927
$ws4.Range("F39").Value = 204
